# The document has a "first page" header/footer that differs from the
# "default" (all other pages) header/footer. Three logo pictures need to be
# renamed (their OOXML docPr/cNvPr "name" attribute only - the images
# themselves, their relationship ids and the rest of the formatting are
# untouched):
#   - First-page header picture (BTec logo, .jpg): image1.jpg -> image2.jpg
#   - Default footer picture   (Pearson logo, .png): image2.png -> image1.png
#   - First-page footer picture (Pearson logo, .png): image2.png -> image1.png
#
# wdHeaderFooterPrimary   = 1
# wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# --- Default (primary) footer: Pearson logo, image2.png -> image1.png ---
$defaultFooter = $section.Footers.Item(1)
if ($defaultFooter.Range.InlineShapes.Count -ge 1) {
    $defaultFooterLogo = $defaultFooter.Range.InlineShapes.Item(1)
    $defaultFooterLogo.Name = "image1.png"
}

# --- First-page footer: Pearson logo, image2.png -> image1.png ---
$section = $d.Sections.Item(1)
$firstFooter = $section.Footers.Item(2)
if ($firstFooter.Range.InlineShapes.Count -ge 1) {
    $firstFooterLogo = $firstFooter.Range.InlineShapes.Item(1)
    $firstFooterLogo.Name = "image1.png"
}

# --- First-page header: BTEC logo, image1.jpg -> image2.jpg ---
$section = $d.Sections.Item(1)
$firstHeader = $section.Headers.Item(2)
if ($firstHeader.Range.InlineShapes.Count -ge 1) {
    $firstHeaderLogo = $firstHeader.Range.InlineShapes.Item(1)
    $firstHeaderLogo.Name = "image2.jpg"
}
